{"js": "// Portuguese (pt) translation of the partner \"reminder to submit documents\" email template.\n// Applies the same text substitutions as the source diff, run by run / occurrence by\n// occurrence, while leaving formatting (highlights, bold, hyperlinks, etc.) untouched.\n\nasync function replaceOccurrence(context, searchText, index, newText, options) {\n  // Re-run the search fresh (indices are stable/document-ordered; doing a fresh\n  // search before each replacement keeps this robust even though we don't\n  // strictly need to since text lengths differ per-call and we always address\n  // by position in the still-untouched parts of the document).\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const body = context.document.body;\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length <= index) {\n    throw new Error(\n      \"replaceOccurrence: expected occurrence \" + index + \" of \" + JSON.stringify(searchText) +\n      \" but only found \" + results.items.length\n    );\n  }\n  results.items[index].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function replaceAll(context, searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const body = context.document.body;\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Standalone \"English\" heading (not the hyperlinked one at the very top).\nawait replaceOccurrence(context, \"English\", 1, \"Ingl\u00eas\");\n\n// 2. \"Subject line\" label \u2014 appears twice, both become \"Linha de assunto\".\nawait replaceAll(context, \"Subject line\", \"Linha de assunto\");\n\n// 3. \"[EVENT NAME]\" placeholder \u2014 appears 4 times, all become \"[NOME DO EVENTO]\".\nawait replaceAll(context, \"[EVENT NAME]\", \"[NOME DO EVENTO]\");\n\n// 4. Subject line tail \u2014 appears twice, both translated.\nawait replaceAll(\n  context,\n  \" \\u2014 have you submitted your docs?  \",\n  \" - j\u00e1 enviou os seus documentos?  \"\n);\n\n// 5. \"Don't forget to send your documents\" heading \u2014 appears twice.\nawait replaceAll(\n  context,\n  \"Don\\u2019t forget to send your documents\",\n  \"N\u00e3o se esque\u00e7a de enviar os seus documentos\"\n);\n\n// 6. \"Hi \" greeting \u2014 only in the first email variant.\nawait replaceOccurrence(context, \"Hi \", 0, \"Ol\u00e1 \");\n\n// 7. \"[PARTNER NAME]\" placeholder \u2014 appears twice.\nawait replaceAll(context, \"[PARTNER NAME]\", \"[NOME DO PARCEIRO]\");\n\n// 8. \"Dear \" greeting \u2014 only in the second email variant.\nawait replaceOccurrence(context, \"Dear \", 0, \"Ol\u00e1 \");\n\n// 9. \"We're excited to see you at the upcoming \" \u2014 appears twice.\nawait replaceAll(\n  context,\n  \"We\\u2019re excited to see you at the upcoming \",\n  \"Estamos ansiosos por v\u00ea-lo no pr\u00f3ximo evento \"\n);\n\n// 10. First-variant-only sentence about confirming registration.\nawait replaceOccurrence(\n  context,\n  \"To confirm your registration, we need the following documents from you by \",\n  0,\n  \"Para confirmar a sua inscri\u00e7\u00e3o, precisamos que nos envie os seguintes documentos at\u00e9 dia \"\n);\n\n// 11. \"[insert list of documents required]\" \u2014 appears twice, and a period is appended.\nawait replaceAll(\n  context,\n  \"[insert list of documents required]\",\n  \"[inserir lista dos documentos necess\u00e1rios].\"\n);\n\n// 12. First-variant-only \"Please send a copy...\" sentence start.\nawait replaceOccurrence(\n  context,\n  \"Please send a copy of these documents to your country manager, \",\n  0,\n  \"Por favor, envie uma c\u00f3pia destes documentos ao seu gestor de parcerias, \"\n);\n\n// 13. \", at \" \u2014 only the first occurrence (first variant) changes.\nawait replaceOccurrence(context, \", at \", 0, \", para \");\n\n// 14. \" or \" \u2014 two of the three occurrences change (1st-variant email + 2nd-variant hyperlink line).\n// Each call re-searches the (already partially edited) document fresh, so once the\n// first match has been turned into \" ou \" it drops out of subsequent result sets \u2014\n// meaning the *next* one to fix is again at index 0.\nawait replaceOccurrence(context, \" or \", 0, \" ou \");\nawait replaceOccurrence(context, \" or \", 0, \" ou \");\n\n// 15. First-variant-only WhatsApp/arrangements sentence tail.\nawait replaceOccurrence(\n  context,\n  \" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\",\n  0,\n  \" (WhatsApp), para podermos tomar as medidas necess\u00e1rias, incluindo alojamento e transporte.\"\n);\n\n// 16. First-variant-only \"If you have any questions, please contact your country manager.\" (no trailing comma/name).\nawait replaceOccurrence(\n  context,\n  \"If you have any questions, please contact your country manager.\",\n  0,\n  \"Se tiver alguma d\u00favida, contacte o gestor do seu pa\u00eds.\"\n);\n\n// 17. \"We look forward to seeing you there!\" \u2014 appears twice.\nawait replaceAll(\n  context,\n  \"We look forward to seeing you there!\",\n  \"Esperamos v\u00ea-lo em breve!\"\n);\n\n// 18. Second-variant-only sentence about the best experience.\nawait replaceOccurrence(\n  context,\n  \"To ensure you have the best experience at this event, we need the following documents from you by \",\n  0,\n  \"De forma a garantir a melhor experi\u00eancia poss\u00edvel neste evento, \u00e9 necess\u00e1rio que nos envie os seguintes documentos at\u00e9 \"\n);\n\n// 19. \"DD Mmm YYYY\" \u2014 only the second variant's date placeholder changes.\nawait replaceOccurrence(context, \"DD Mmm YYYY\", 1, \"DD Mmm AAAA\");\n\n// 20. Second-variant-only \"Please reply to this email...\" sentence.\nawait replaceOccurrence(\n  context,\n  \"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\",\n  0,\n  \"Por favor, responda a este e-mail com uma c\u00f3pia destes documentos para que possamos tomar as provid\u00eancias necess\u00e1rias, incluindo alojamento e transporte.\"\n);\n\n// 21. Second-variant-only \"If you have any questions, please contact us via \" (before the hyperlinks).\nawait replaceOccurrence(\n  context,\n  \"If you have any questions, please contact us via \",\n  0,\n  \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \"\n);\n", "ps1": "# Portuguese (pt) translation of the partner \"reminder to submit documents\" email\n# template. Applies the same text substitutions as the source diff, scoped to\n# specific occurrences (by position in the document) so formatting/hyperlinks/\n# comments are left untouched and only the intended runs are translated.\n\n$d = $word.ActiveDocument\n\nfunction Replace-NthOccurrence($SearchText, $ReplaceText, $Occurrence) {\n    # 0-based occurrence index, scanning the whole document from the top.\n    $r = $d.Content\n    $count = 0\n    while ($true) {\n        $found = $r.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n        if (-not $found) {\n            throw \"Replace-NthOccurrence: could not find occurrence $Occurrence of '$SearchText' (only saw $count)\"\n        }\n        if ($count -eq $Occurrence) {\n            $r.Text = $ReplaceText\n            return\n        }\n        $count = $count + 1\n    }\n}\n\nfunction Replace-AllOccurrences($SearchText, $ReplaceText) {\n    $r = $d.Content\n    while ($true) {\n        $found = $r.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n        if (-not $found) {\n            return\n        }\n        $r.Text = $ReplaceText\n    }\n}\n\n# 1. Standalone \"English\" heading (not the hyperlinked one at the very top).\nReplace-NthOccurrence \"English\" \"Ingl\u00eas\" 1\n\n# 2. \"Subject line\" label \u2014 appears twice, both become \"Linha de assunto\".\nReplace-AllOccurrences \"Subject line\" \"Linha de assunto\"\n\n# 3. \"[EVENT NAME]\" placeholder \u2014 appears 4 times, all become \"[NOME DO EVENTO]\".\nReplace-AllOccurrences \"[EVENT NAME]\" \"[NOME DO EVENTO]\"\n\n# 4. Subject line tail \u2014 appears twice, both translated.\nReplace-AllOccurrences \" \u2014 have you submitted your docs?  \" \" - j\u00e1 enviou os seus documentos?  \"\n\n# 5. \"Don't forget to send your documents\" heading \u2014 appears twice.\nReplace-AllOccurrences \"Don\u2019t forget to send your documents\" \"N\u00e3o se esque\u00e7a de enviar os seus documentos\"\n\n# 6. \"Hi \" greeting \u2014 only in the first email variant.\nReplace-NthOccurrence \"Hi \" \"Ol\u00e1 \" 0\n\n# 7. \"[PARTNER NAME]\" placeholder \u2014 appears twice.\nReplace-AllOccurrences \"[PARTNER NAME]\" \"[NOME DO PARCEIRO]\"\n\n# 8. \"Dear \" greeting \u2014 only in the second email variant.\nReplace-NthOccurrence \"Dear \" \"Ol\u00e1 \" 0\n\n# 9. \"We're excited to see you at the upcoming \" \u2014 appears twice.\nReplace-AllOccurrences \"We\u2019re excited to see you at the upcoming \" \"Estamos ansiosos por v\u00ea-lo no pr\u00f3ximo evento \"\n\n# 10. First-variant-only sentence about confirming registration.\nReplace-NthOccurrence \"To confirm your registration, we need the following documents from you by \" \"Para confirmar a sua inscri\u00e7\u00e3o, precisamos que nos envie os seguintes documentos at\u00e9 dia \" 0\n\n# 11. \"[insert list of documents required]\" \u2014 appears twice, and a period is appended.\nReplace-AllOccurrences \"[insert list of documents required]\" \"[inserir lista dos documentos necess\u00e1rios].\"\n\n# 12. First-variant-only \"Please send a copy...\" sentence start.\nReplace-NthOccurrence \"Please send a copy of these documents to your country manager, \" \"Por favor, envie uma c\u00f3pia destes documentos ao seu gestor de parcerias, \" 0\n\n# 13. \", at \" \u2014 only the first occurrence (first variant) changes.\nReplace-NthOccurrence \", at \" \", para \" 0\n\n# 14. \" or \" \u2014 two of the three occurrences change (1st-variant email + 2nd-variant hyperlink line).\n# Each call re-scans the (already partially edited) document fresh, so once the\n# first match has been turned into \" ou \" it drops out of the scan \u2014 meaning the\n# *next* one to fix is again at index 0.\nReplace-NthOccurrence \" or \" \" ou \" 0\nReplace-NthOccurrence \" or \" \" ou \" 0\n\n# 15. First-variant-only WhatsApp/arrangements sentence tail.\nReplace-NthOccurrence \" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\" \" (WhatsApp), para podermos tomar as medidas necess\u00e1rias, incluindo alojamento e transporte.\" 0\n\n# 16. First-variant-only \"If you have any questions, please contact your country manager.\" (no trailing comma/name).\nReplace-NthOccurrence \"If you have any questions, please contact your country manager.\" \"Se tiver alguma d\u00favida, contacte o gestor do seu pa\u00eds.\" 0\n\n# 17. \"We look forward to seeing you there!\" \u2014 appears twice.\nReplace-AllOccurrences \"We look forward to seeing you there!\" \"Esperamos v\u00ea-lo em breve!\"\n\n# 18. Second-variant-only sentence about the best experience.\nReplace-NthOccurrence \"To ensure you have the best experience at this event, we need the following documents from you by \" \"De forma a garantir a melhor experi\u00eancia poss\u00edvel neste evento, \u00e9 necess\u00e1rio que nos envie os seguintes documentos at\u00e9 \" 0\n\n# 19. \"DD Mmm YYYY\" \u2014 only the second variant's date placeholder changes.\nReplace-NthOccurrence \"DD Mmm YYYY\" \"DD Mmm AAAA\" 1\n\n# 20. Second-variant-only \"Please reply to this email...\" sentence.\nReplace-NthOccurrence \"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\" \"Por favor, responda a este e-mail com uma c\u00f3pia destes documentos para que possamos tomar as provid\u00eancias necess\u00e1rias, incluindo alojamento e transporte.\" 0\n\n# 21. Second-variant-only \"If you have any questions, please contact us via \" (before the hyperlinks).\nReplace-NthOccurrence \"If you have any questions, please contact us via \" \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \" 0\n"}
